$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "GoalTime" (row 19) and "GoalDistance" (row 20) entries which are
# unused variables. Selecting the two entire rows and deleting them shifts
# everything below up by two rows (matches the diff's row renumbering).
$rows = $ws.Range("A19:XFD20")
$rows.Select()
$rows.Delete()
